$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sarah")

# Add the new order row (row 5)
$ws.Range("A5").Value = "Sarah"
$ws.Range("B5").Value = "Sarah_Email"
$ws.Range("C5").Value = "Northwoods Cranberry Sauce"
$ws.Range("D5").Value = 6

# Widen column C to fit the new, longer item name
$ws.Columns.Item(3).ColumnWidth = 15.85

# Move the active selection (as left by the user after entering data)
$ws.Range("E8").Select()
